$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 1.02
    "C2" = 1.039369959750353
    "D2" = 1.045329498178674
    "E2" = 1.047346641722157
    "F2" = 1.05685754876978
    "I2" = 1.037451470814697
    "J2" = 1.044462345111605
    "K2" = 1.048098043204668
    "L2" = 1.050109535099643
    "M2" = 1.059594109347405
    "N2" = 1.005712725503983
    "B3" = 1.02
    "C3" = 1.040298039018768
    "D3" = 1.046016249075007
    "E3" = 1.048144585214739
    "F3" = 1.057674394032931
    "I3" = 1.037593281652592
    "J3" = 1.045035709184041
    "K3" = 1.048596537367497
    "L3" = 1.05071934034515
    "M3" = 1.060224668439941
    "B4" = 1.02
    "C4" = 1.040899290802825
    "D4" = 1.046461067549612
    "E4" = 1.048661859337786
    "F4" = 1.058203792220208
    "I4" = 1.037683906595885
    "J4" = 1.045406803599949
    "K4" = 1.048918861989701
    "L4" = 1.051114219638495
    "M4" = 1.06063287288626
    "B5" = 1.02
    "C5" = 1.041152228640804
    "D5" = 1.046648174014057
    "E5" = 1.048879547430178
    "F5" = 1.058426551644131
    "I5" = 1.037721732689797
    "J5" = 1.045562831886991
    "K5" = 1.049054309870928
    "L5" = 1.051280296052598
    "M5" = 1.06080452618454
    "B6" = 1.02
    "C6" = 1.041194708039193
    "D6" = 1.046679596107805
    "E6" = 1.048916111425125
    "F6" = 1.058463965632557
    "I6" = 1.037728067851678
    "J6" = 1.045589030871219
    "K6" = 1.049077048764012
    "L6" = 1.051308185029576
    "M6" = 1.060833350074855
    "B7" = 1.02
    "C7" = 1.040902669898075
    "D7" = 1.04646356726577
    "E7" = 1.04866476721064
    "F7" = 1.058206767957249
    "I7" = 1.037684413102054
    "J7" = 1.045408888379759
    "K7" = 1.048920672078013
    "L7" = 1.051116438489616
    "M7" = 1.060635166354403
    "B8" = 1.02
    "C8" = 1.03968345807647
    "D8" = 1.045561495761002
    "E8" = 1.047616112461955
    "F8" = 1.057133429476542
    "I8" = 1.037499631316672
    "J8" = 1.044656096904438
    "K8" = 1.04826655942103
    "L8" = 1.050315559842656
    "M8" = 1.059807169407396
    "B9" = 1.02
    "C9" = 1.037540644539926
    "D9" = 1.043975411225252
    "E9" = 1.045775609913113
    "F9" = 1.055248618418613
    "I9" = 1.037165348757322
    "J9" = 1.043330322460074
    "K9" = 1.047112189810763
    "L9" = 1.048906624723262
    "M9" = 1.058349657305174
    "B10" = 1.02
    "C10" = 1.036115943787896
    "D10" = 1.042920462350671
    "E10" = 1.044553658494293
    "F10" = 1.053996587244456
    "I10" = 1.036936700029401
    "J10" = 1.042447044078302
    "K10" = 1.046341515791821
    "L10" = 1.047968971671646
    "M10" = 1.057379094814896
    "B11" = 1.02
    "C11" = 1.03549996065734
    "D11" = 1.042464257962953
    "E11" = 1.044025759185668
    "F11" = 1.053455535995426
    "I11" = 1.036836326363778
    "J11" = 1.042064725702049
    "K11" = 1.046007561708655
    "L11" = 1.047563362707433
    "M11" = 1.056959112730128
    "B12" = 1.02
    "C12" = 1.035271296901633
    "D12" = 1.042294894506254
    "E12" = 1.043829857943281
    "F12" = 1.053254730519477
    "I12" = 1.036798838320858
    "J12" = 1.041922738948623
    "K12" = 1.045883480424825
    "L12" = 1.047412763076285
    "M12" = 1.056803155905701
    "B13" = 1.02
    "C13" = 1.035320339705525
    "D13" = 1.042331219389806
    "E13" = 1.043871871067845
    "F13" = 1.053297796477848
    "I13" = 1.036806888890718
    "J13" = 1.041953194514436
    "K13" = 1.045910097884234
    "L13" = 1.047445064399229
    "M13" = 1.056836607194424
    "B14" = 1.02
    "C14" = 1.03548105638025
    "D14" = 1.042450256457172
    "E14" = 1.044009562163046
    "F14" = 1.053438933974731
    "I14" = 1.036833231764603
    "J14" = 1.042052988549995
    "K14" = 1.045997305835418
    "L14" = 1.047550912830539
    "M14" = 1.056946220400141
    "B15" = 1.02
    "C15" = 1.035580097793383
    "D15" = 1.042523611254452
    "E15" = 1.044094422633029
    "F15" = 1.053525915369376
    "I15" = 1.036849435359978
    "J15" = 1.042114478075105
    "K15" = 1.04605103281729
    "L15" = 1.047616137745302
    "M15" = 1.057013762485229
    "B16" = 1.02
    "C16" = 1.036156844092171
    "D16" = 1.042950751824388
    "E16" = 1.044588719151139
    "F16" = 1.054032518080192
    "I16" = 1.036943332745412
    "J16" = 1.04247242049175
    "K16" = 1.046363674125168
    "L16" = 1.047995899166741
    "M16" = 1.057406973623034
    "B17" = 1.02
    "C17" = 1.036518869510398
    "D17" = 1.043218846440379
    "E17" = 1.044899104457876
    "F17" = 1.054350588745859
    "I17" = 1.037001866478546
    "J17" = 1.042696988535167
    "K17" = 1.046559720480956
    "L17" = 1.048234221785966
    "M17" = 1.057653700017349
    "B18" = 1.02
    "C18" = 1.036730121652672
    "D18" = 1.043375278830633
    "E18" = 1.045080263947129
    "F18" = 1.054536218694915
    "I18" = 1.037035876252537
    "J18" = 1.04282798937551
    "K18" = 1.046674047072593
    "L18" = 1.048373269937082
    "M18" = 1.057797638061267
    "B19" = 1.02
    "C19" = 1.036802168202467
    "D19" = 1.043428627906431
    "E19" = 1.045142054442839
    "F19" = 1.054599531416117
    "I19" = 1.037047450304634
    "J19" = 1.042872659630892
    "K19" = 1.046713025388371
    "L19" = 1.048420688269262
    "M19" = 1.057846721700613
    "B20" = 1.02
    "C20" = 1.036480018410351
    "D20" = 1.043190076482385
    "E20" = 1.044865790938573
    "F20" = 1.054316451929777
    "I20" = 1.036995600007931
    "J20" = 1.04267289304208
    "K20" = 1.046538689028501
    "L20" = 1.048208648022093
    "M20" = 1.057627225836604
    "B21" = 1.02
    "C21" = 1.035433725455264
    "D21" = 1.042415200477434
    "E21" = 1.043969010472464
    "F21" = 1.053397367928967
    "I21" = 1.036825480093159
    "J21" = 1.04202360104035
    "K21" = 1.045971626247383
    "L21" = 1.047519741394449
    "M21" = 1.056913940838133
    "B22" = 1.02
    "C22" = 1.034776689082257
    "D22" = 1.041928533477486
    "E22" = 1.043406234796854
    "F22" = 1.052820459154248
    "I22" = 1.036717334381913
    "J22" = 1.041615500887315
    "K22" = 1.045614884474229
    "L22" = 1.047086955511406
    "M22" = 1.056465721045479
    "B23" = 1.02
    "C23" = 1.035124919314671
    "D23" = 1.04218647421572
    "E23" = 1.04370447123186
    "F23" = 1.053126198206019
    "I23" = 1.036774776545445
    "J23" = 1.041831829245137
    "K23" = 1.045804019227504
    "L23" = 1.047316349199926
    "M23" = 1.056703306573549
    "B24" = 1.02
    "C24" = 1.036497573275765
    "D24" = 1.043203076211062
    "E24" = 1.044880843522708
    "F24" = 1.054331876564242
    "I24" = 1.036998431964077
    "J24" = 1.042683780714799
    "K24" = 1.046548192310665
    "L24" = 1.048220203586069
    "M24" = 1.057639188296951
    "B25" = 1.02
    "C25" = 1.038093942032549
    "D25" = 1.04438502924443
    "E25" = 1.046250541192182
    "F25" = 1.055735100230883
    "I25" = 1.037252792765049
    "J25" = 1.043672971631012
    "K25" = 1.047410820412209
    "L25" = 1.049270585566157
    "M25" = 1.058726269477798
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

Write-Host "Updated $($values.Count) cells"